$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 689.59375
$ws.Range("J121").Value = 662.4
$ws.Range("L121").Value = 1987.2
$ws.Range("N121").Value = -5481.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 961.5098
$ws.Range("J129").Value = 1129.7632
$ws.Range("L129").Value = 3389.2896
$ws.Range("N129").Value = -13389.2896

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2423.2327
$ws.Range("I132").Value = 2153.2163
$ws.Range("J132").Value = 4088.3333
$ws.Range("K132").Value = 6459.6489
$ws.Range("L132").Value = 12264.9999
$ws.Range("M132").Value = -3929.6489
$ws.Range("N132").Value = -17324.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16680.957
$ws.Range("I32").Value = 13226.556
$ws.Range("J32").Value = 47770.57
$ws.Range("K32").Value = 13226.556
$ws.Range("L32").Value = 47770.57
$ws.Range("M32").Value = -12939.556
$ws.Range("N32").Value = -48344.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1963
$ws.Range("I61").Value = 1535.4546
$ws.Range("K61").Value = 1535.4546
$ws.Range("M61").Value = -1323.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1074.2826
$ws.Range("I74").Value = 909.7179599999999
$ws.Range("J74").Value = 1991.1428
$ws.Range("K74").Value = 909.7179599999999
$ws.Range("L74").Value = 1991.1428
$ws.Range("M74").Value = -35.71795999999995
$ws.Range("N74").Value = -3739.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1074.2826
$ws.Range("I77").Value = 909.7179599999999
$ws.Range("J77").Value = 1991.1428
$ws.Range("K77").Value = 4548.5898
$ws.Range("L77").Value = 9955.714
$ws.Range("M77").Value = -180.5897999999997
$ws.Range("N77").Value = -18691.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 535.9091
$ws.Range("I97").Value = 535.9091
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 535.9091
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -39.90909999999997
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2250.8333
$ws.Range("I122").Value = 2307.1875
$ws.Range("K122").Value = 6921.5625
$ws.Range("M122").Value = -4471.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 529864
$ws.Range("I132").Value = 556820
$ws.Range("K132").Value = 1670460
$ws.Range("M132").Value = -1667930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 49820.9
$ws.Range("J134").Value = 49820.9
$ws.Range("L134").Value = 49820.9
$ws.Range("N134").Value = -59960.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1963
$ws.Range("I136").Value = 1535.4546
$ws.Range("K136").Value = 4606.3638
$ws.Range("M136").Value = -2056.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2178880.8
$ws.Range("I80").Value = 5291210.5
$ws.Range("J80").Value = 250
$ws.Range("K80").Value = 5291210.5
$ws.Range("L80").Value = 250
$ws.Range("M80").Value = -5290212.5
$ws.Range("N80").Value = -2246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 2178880.8
$ws.Range("I83").Value = 5291210.5
$ws.Range("J83").Value = 250
$ws.Range("K83").Value = 26456052.5
$ws.Range("L83").Value = 1250
$ws.Range("M83").Value = -26451060.5
$ws.Range("N83").Value = -11234

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1550.7273
$ws.Range("I99").Value = 843.375
$ws.Range("J99").Value = 3437
$ws.Range("K99").Value = 843.375
$ws.Range("L99").Value = 3437
$ws.Range("M99").Value = 654.625
$ws.Range("N99").Value = -6433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 54554.74
$ws.Range("J130").Value = 54554.74
$ws.Range("L130").Value = 54554.74
$ws.Range("N130").Value = -64594.74

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 346977.38
$ws.Range("I134").Value = 436352.1
$ws.Range("J134").Value = 4374.1665
$ws.Range("K134").Value = 1309056.3
$ws.Range("L134").Value = 13122.4995
$ws.Range("M134").Value = -1306521.3
$ws.Range("N134").Value = -18192.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8944.111000000001
$ws.Range("J4").Value = 9999.625
$ws.Range("L4").Value = 9999.625
$ws.Range("N4").Value = -10223.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 15000
$ws.Range("J108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("N108").Value = -27680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100487.4
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 12000
$ws.Range("N4").Value = -12224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3300
$ws.Range("N113").Value = -7640
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1878.2778
$ws.Range("I132").Value = 1219.9
$ws.Range("J132").Value = 2701.25
$ws.Range("K132").Value = 10979.1
$ws.Range("L132").Value = 24311.25
$ws.Range("M132").Value = -8449.1
$ws.Range("N132").Value = -29371.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 24002
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 47801.91
$ws.Range("I97").Value = 73115.71000000001
$ws.Range("J97").Value = 3502.75
$ws.Range("K97").Value = 73115.71000000001
$ws.Range("L97").Value = 3502.75
$ws.Range("M97").Value = -72619.71000000001
$ws.Range("N97").Value = -4494.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2938.3845
$ws.Range("I102").Value = 2969.9
$ws.Range("J102").Value = 2833.3333
$ws.Range("K102").Value = 2969.9
$ws.Range("L102").Value = 2833.3333
$ws.Range("M102").Value = -1347.9
$ws.Range("N102").Value = -6077.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1602.3334
$ws.Range("I113").Value = 984.2
$ws.Range("J113").Value = 2375
$ws.Range("K113").Value = 984.2
$ws.Range("L113").Value = 2375
$ws.Range("M113").Value = 1185.8
$ws.Range("N113").Value = -6715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4995.2
$ws.Range("I122").Value = 4619
$ws.Range("K122").Value = 13857
$ws.Range("M122").Value = -11407

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11844.444
$ws.Range("J123").Value = 11844.444
$ws.Range("L123").Value = 11844.444
$ws.Range("N123").Value = -16744.444

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 41499.332
$ws.Range("J129").Value = 41499.332
$ws.Range("L129").Value = 41499.332
$ws.Range("N129").Value = -51499.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1876.5
$ws.Range("I132").Value = 1352.1
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 4056.3
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -1526.3
$ws.Range("N132").Value = -18555.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9467.5
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 9995
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 9995
$ws.Range("M2").Value = -388
$ws.Range("N2").Value = -10219

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 37500
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5257.5806
$ws.Range("I122").Value = 5607.9546
$ws.Range("J122").Value = 4401.1113
$ws.Range("K122").Value = 16823.8638
$ws.Range("L122").Value = 13203.3339
$ws.Range("M122").Value = -14373.8638
$ws.Range("N122").Value = -18103.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 55076.332
$ws.Range("J129").Value = 55076.332
$ws.Range("L129").Value = 55076.332
$ws.Range("N129").Value = -65076.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 25389.727
$ws.Range("J130").Value = 25389.727
$ws.Range("L130").Value = 25389.727
$ws.Range("N130").Value = -35429.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2151.4348
$ws.Range("I136").Value = 2089.8823
$ws.Range("J136").Value = 2325.8333
$ws.Range("K136").Value = 6269.646900000001
$ws.Range("L136").Value = 6977.499899999999
$ws.Range("M136").Value = -3719.646900000001
$ws.Range("N136").Value = -12077.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3400
$ws.Range("J2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("N2").Value = -3724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25056.666
$ws.Range("J64").Value = 25056.666
$ws.Range("L64").Value = 25056.666
$ws.Range("N64").Value = -25552.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 25056.666
$ws.Range("J67").Value = 25056.666
$ws.Range("L67").Value = 25056.666
$ws.Range("N67").Value = -26772.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 20715
$ws.Range("J128").Value = 20715
$ws.Range("L128").Value = 20715
$ws.Range("N128").Value = -30675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 23000
$ws.Range("J129").Value = 23000
$ws.Range("L129").Value = 23000
$ws.Range("N129").Value = -33000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 42965.6
$ws.Range("J130").Value = 42965.6
$ws.Range("L130").Value = 42965.6
$ws.Range("N130").Value = -53005.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 53231.715
$ws.Range("J131").Value = 53231.715
$ws.Range("L131").Value = 53231.715
$ws.Range("N131").Value = -63311.715

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1689
$ws.Range("I132").Value = 1174.4546
$ws.Range("J132").Value = 3232.6365
$ws.Range("K132").Value = 3523.3638
$ws.Range("L132").Value = 9697.9095
$ws.Range("M132").Value = -993.3638000000001
$ws.Range("N132").Value = -14757.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1240.7675
$ws.Range("I136").Value = 1075.8
$ws.Range("J136").Value = 1962.5
$ws.Range("K136").Value = 3227.4
$ws.Range("L136").Value = 5887.5
$ws.Range("M136").Value = -677.3999999999996
$ws.Range("N136").Value = -10987.5
